$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing scores for row 4 (C4:F4) so J4's SUM formula recalculates to 20
$ws.Range("C4").Value = 5
$ws.Range("D4").Value = 5
$ws.Range("E4").Value = 5
$ws.Range("F4").Value = 5

# Move the active selection to G6 (also brings the frozen pane's
# top-left cell back up to C4 since the sheet scrolls to show it)
$ws.Range("G6").Select()
